# Interdiff between v7 and v8
#
# 1) The cached text of every auto-updating "datetimeFigureOut" date field
#    (slide master + all 11 slide layouts) moves from 2/2/2017 -> 2/5/2017.
# 2) The "TextBox 78" shape on slide 1 moves down slightly
#    (y: 3472934 EMU -> 3657600 EMU, i.e. Top: 273.4594pt -> 288pt).

$p = $ppt.ActivePresentation

$oldDate = "2/2/2017"
$newDate = "2/5/2017"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $sh = $container.Shapes.Item($j)
        if (-not $sh.HasTextFrame) { continue }

        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and ($sh.TextFrame.TextRange.Text -eq $oldDate)) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master

# Every Custom Layout (slideLayout1.xml .. slideLayout11.xml) hanging off
# the slide master.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout
}

# Reposition "TextBox 78" on slide 1.
$slide = $p.Slides.Item(1)
$box = $slide.Shapes.Item("TextBox 78")
$box.Top = 288
